$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 20: header block (bold), mirrors row 1 but for a new Milestone4 run ----
$ws.Range("A20").Value = "Milestone4"
$ws.Range("C20").Value = "gamma = 0.5"
$ws.Range("E20").Value = "alpha = 0.4"
$ws.Range("G20").Value = "256 states"
$ws.Range("H20").Value = "   24576 actions"
$ws.Range("A20:I20").Font.Bold = $true

# ---- Row 21: training bucket descriptions ----
$ws.Range("A21").Value = "0-9000 actions 60% random"
$ws.Range("D21").Value = "9000-18000 actions 40% random"
$ws.Range("H21").Value = "18000-27000 actions 20% random"

# ---- Row 22: counts ----
$ws.Range("H22").Value = "(5011/9000)"
$ws.Range("D22").Value = "  (3193/9000)"
$ws.Range("B22").Value = " (1217/9000)"

# ---- Row 23: Training results ----
$ws.Range("A23").Value = "Training"
$ws.Range("B23").Value = 0.135
$ws.Range("B23").NumberFormat = "0%"
$ws.Range("D23").Value = 0.355
$ws.Range("D23").NumberFormat = "0.00%"
$ws.Range("H23").Value = 0.557
$ws.Range("H23").NumberFormat = "0.00%"

# ---- Row 24: Evaluation result ----
$ws.Range("A24").Value = "Evaluation (0% random)"
$ws.Range("D24").Value = "74.3% (1840/2477)  6000 actions "

# ---- Final selection ----
$ws.Range("F31").Select()
